$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 3, shifting existing rows 3-19 down to 4-20
$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = 4
$ws.Range("B3").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C3").Value = "Los Lagos"
$ws.Range("D3").Value = 44490
$ws.Range("E3").Value = 10
$ws.Range("F3").Value = 300000000
$ws.Range("G3").Value = "Espárragos"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 72
$ws.Range("K3").Value = 1700
$ws.Range("L3").Value = 1700
$ws.Range("M3").Value = 1700
$ws.Range("N3").Value = "$/kilo"
$ws.Range("O3").Value = "Provincia de Linares"
$ws.Range("P3").Value = 1700
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = "Hortaliza"
